$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 5 (pushing the previously
# existing rows 5-104 down to rows 6-105). Insert a fresh row so all the
# formatting/styles of the existing row 5 (notably the date style on column D)
# shift down along with the data.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new record's data.
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(5, 4).Value = 44963
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = 100112040
$ws.Cells.Item(5, 7).Value = "Cilantro"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 300
$ws.Cells.Item(5, 11).Value = 2000
$ws.Cells.Item(5, 12).Value = 2500
$ws.Cells.Item(5, 13).Value = 2250
$ws.Cells.Item(5, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(5, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5, 16).Value = 1125
$ws.Cells.Item(5, 17).Value = 2
$ws.Cells.Item(5, 18).Value = "Hortaliza"
